# Update the "May" sales data: rewrite the Amount (column E) figures for
# rows 7-63 so the bar chart's values/labels line up, and move the thick
# "group divider" bottom-border down by one row in each of the affected
# five-row blocks (it previously sat one row too high relative to the
# Date/category groups).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New Amount values for E7:E63
$newValues = @{
    7  = 12000;  8  = 15000;  9  = 20000;  10 = 25000;  11 = 30000;
    12 = 40000;  13 = 22000;  14 = 26000;  15 = 18000;  16 = 17000;
    17 = 23000;  18 = 14000;  19 = 22000;  20 = 26000;  21 = 18000;
    22 = 17000;  23 = 23000;  24 = 14000;  25 = 22000;  26 = 26000;
    27 = 18000;  28 = 17000;  29 = 23000;  30 = 14000;  31 = 22000;
    32 = 26000;  33 = 18000;  34 = 17000;  35 = 23000;  36 = 14000;
    37 = 22000;  38 = 26000;  39 = 18000;  40 = 17000;  41 = 23000;
    42 = 14000;  43 = 22000;  44 = 26000;  45 = 18000;  46 = 17000;
    47 = 23000;  48 = 14000;  49 = 22000;  50 = 26000;  51 = 18000;
    52 = 17000;  53 = 23000;  54 = 14000;  55 = 22000;  56 = 26000;
    57 = 18000;  58 = 17000;  59 = 23000;  60 = 14000;  61 = 15000;
    62 = 13500;  63 = 12000
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 5).Value2 = $newValues[$row]
}

# Rows whose E-cell needs the thick (medium) bottom border turned ON
# (the group divider moves onto this row).
$addThickBottom = @(11, 24, 30, 36, 42, 54, 60)
foreach ($row in $addThickBottom) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = -4138
}

# Rows whose E-cell needs the thick bottom border turned back OFF
# (reverts to the regular thin grid border used elsewhere in column E).
$removeThickBottom = @(23, 28, 33, 38, 43, 53, 58)
foreach ($row in $removeThickBottom) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = 2
}

# Match the saved selection/active-cell state from the edit.
$ws.Range("E58:E63").Select() | Out-Null
